# "sume todas las variables de Wunder"
# Insert a new sheet ("Sheet1") right after "estaciones" and fill it with
# every station code (Wunder variable) that appears in column B of
# "estaciones" - i.e. "add up all the Wunder variables" into their own
# sheet. Also nudges the remembered selection on a few sheets, matching
# the author's saved cursor positions.

$wb = $excel.ActiveWorkbook

$estaciones = $wb.Worksheets.Item("estaciones")

# 1) New sheet, placed immediately after "estaciones".
$nuevo = $wb.Worksheets.Add($null, $estaciones)

# Worksheet handles resolve by (then-current) position, so fetch the
# remaining sheets by name only AFTER the insert has shifted everyone
# else's index - otherwise $db/$apiKeys/$tamanoDB would silently drift to
# whatever sheet now occupies their old slot.
$db         = $wb.Worksheets.Item("db")
$apiKeys    = $wb.Worksheets.Item("apiKeys")
$tamanoDB   = $wb.Worksheets.Item("tamaño_DB")

# 2) Copy the 143 station codes from estaciones!B2:B144 into Sheet1!A1:A143.
$estaciones.Range("B2:B144").Copy()
$nuevo.Range("A1").PasteSpecial(-4163)   # xlPasteValues
$excel.CutCopyMode = $false

# 3) Restore/update the remembered selections on each sheet.
$estaciones.Activate()
$estaciones.Range("B28").Select()

$nuevo.Activate()
$nuevo.Range("A130").Select()

$apiKeys.Activate()
$apiKeys.Range("D9").Select()

$tamanoDB.Activate()
$tamanoDB.Range("E31").Select()

$estaciones.Activate()
